# Restored from revision #f9be1eab79f5ed1f98d78ea7ff9c7e9e071fb0b7.TEST Author: admin. Type: SAVE.
# Update cell C10 on the active sheet from 18 to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C10").Value = 1
